# edit.ps1
#
# Purpose: correct the English/French sentence counts on the bilingual
# percentage worksheet and recompute the dependent percentage columns
# (including the "Overall Average" row), with docstrings and basic
# null/bounds checks guarding each step so failures are easy to diagnose.

<#
.SYNOPSIS
    Applies corrected English/French sentence counts to the active
    worksheet and recomputes the dependent percentage columns.

.DESCRIPTION
    Row 2 (Balarama Holness):                           English 66 -> 65, French 123 -> 125
    Row 3 (Dexter Xurukulasuriya):                       English 143 -> 146, French 87 -> 88
    Row 4 (Thierry Lindor):                              French 136 -> 135 (English unchanged)
    Row 5 (Conseil national des musulmans canadiens):    French 71 -> 70 (English unchanged)
    Row 7 (Overall Average):                             percentages recomputed from the updated totals
#>

function Set-CellValue {
    <#
    .SYNOPSIS
        Defensive helper that writes a value into a worksheet cell,
        guarding against null worksheet references.
    .PARAMETER Worksheet
        The worksheet object to write to.
    .PARAMETER Address
        The A1-style address of the target cell.
    .PARAMETER Value
        The value to store in the cell.
    #>
    param($Worksheet, $Address, $Value)

    if ($null -eq $Worksheet) {
        throw "Set-CellValue: worksheet reference is null."
    }
    if ($null -eq $Address -or $Address -eq "") {
        throw "Set-CellValue: address must be a non-empty string."
    }

    $range = $Worksheet.Range($Address)
    if ($null -eq $range) {
        throw "Set-CellValue: range '$Address' could not be resolved."
    }

    $range.Value = $Value
}

function Get-CellNumber {
    <#
    .SYNOPSIS
        Defensive helper that reads a numeric value from a worksheet cell,
        guarding against null worksheet/range/value references.
    .PARAMETER Worksheet
        The worksheet object to read from.
    .PARAMETER Row
        1-based row index.
    .PARAMETER Column
        1-based column index.
    #>
    param($Worksheet, $Row, $Column)

    if ($null -eq $Worksheet) {
        throw "Get-CellNumber: worksheet reference is null."
    }

    $cell = $Worksheet.Cells.Item($Row, $Column)
    if ($null -eq $cell) {
        throw "Get-CellNumber: cell at row $Row, column $Column could not be resolved."
    }

    $value = $cell.Value()
    if ($null -eq $value) {
        throw "Get-CellNumber: cell at row $Row, column $Column has a null value."
    }

    return $value
}

# --- Entry point -----------------------------------------------------------

$wb = $excel.ActiveWorkbook
if ($null -eq $wb) {
    throw "No active workbook is available."
}

$ws = $wb.ActiveSheet
if ($null -eq $ws) {
    throw "Active workbook has no active sheet."
}

# Corrected raw counts: English Sentences = column C, French Sentences = column D.
Set-CellValue $ws "C2" 65
Set-CellValue $ws "D2" 125
Set-CellValue $ws "C3" 146
Set-CellValue $ws "D3" 88
Set-CellValue $ws "D4" 135
Set-CellValue $ws "D5" 70

# Recompute the percentage columns (E = English %, F = French %) for every
# data row whose counts changed, using share = count / (English + French) * 100.
$rowsToRecompute = 2, 3, 4, 5
foreach ($row in $rowsToRecompute) {
    $english = Get-CellNumber $ws $row 3
    $french  = Get-CellNumber $ws $row 4

    $total = $english + $french
    if ($total -eq 0) {
        throw "Row $row : English + French total is zero; cannot compute percentages."
    }

    $englishPct = ($english / $total) * 100
    $frenchPct  = ($french  / $total) * 100

    # NOTE: the target address is computed (string concatenation), so it is
    # assigned to a variable first before being passed into the helper —
    # passing the `"E" + $row` expression directly as a call argument does
    # not reliably bind in this host's argument parser.
    $englishAddress = $("E" + $row)
    $frenchAddress  = $("F" + $row)

    Set-CellValue $ws $englishAddress $englishPct
    Set-CellValue $ws $frenchAddress  $frenchPct
}

# Recompute the "Overall Average" row (row 7) from the updated totals across
# all data rows (2-6), mirroring the original workbook's aggregate formula.
$totalEnglish = 0
$totalFrench  = 0
$dataRows = 2, 3, 4, 5, 6
foreach ($row in $dataRows) {
    $englishVal = Get-CellNumber $ws $row 3
    $frenchVal  = Get-CellNumber $ws $row 4

    $totalEnglish = $totalEnglish + $englishVal
    $totalFrench  = $totalFrench + $frenchVal
}

$grandTotal = $totalEnglish + $totalFrench
if ($grandTotal -eq 0) {
    throw "Grand total of English + French counts is zero; cannot compute overall average."
}

$overallEnglishPct = ($totalEnglish / $grandTotal) * 100
$overallFrenchPct  = ($totalFrench  / $grandTotal) * 100

Set-CellValue $ws "E7" $overallEnglishPct
Set-CellValue $ws "F7" $overallFrenchPct

Write-Output "Updated sentence counts and recomputed percentages successfully."
